$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 150
$ws.Range("M4").Value = -36

$ws.Range("H17").Value = 443075.62
$ws.Range("J17").Value = 443075.62
$ws.Range("L17").Value = 1329226.86
$ws.Range("N17").Value = -1329562.86

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H39").Value = 143.38095
$ws.Range("I39").Value = 17.4
$ws.Range("J39").Value = 458.33334
$ws.Range("K39").Value = 52.2
$ws.Range("L39").Value = 1375.00002
$ws.Range("M39").Value = 243.8
$ws.Range("N39").Value = -1967.00002

$ws.Range("H48").Value = 1965
$ws.Range("I48").Value = 1965
$ws.Range("K48").Value = 5895
$ws.Range("M48").Value = -5603

$ws.Range("H51").Value = 2668.182
$ws.Range("I51").Value = 2538.889
$ws.Range("K51").Value = 2538.889
$ws.Range("M51").Value = -2054.889

$ws.Range("H56").Value = 1965
$ws.Range("I56").Value = 1965
$ws.Range("K56").Value = 5895
$ws.Range("M56").Value = -5361

$ws.Range("H58").Value = 1828.5
$ws.Range("I58").Value = 1471.3334
$ws.Range("J58").Value = 2900
$ws.Range("K58").Value = 4414.0002
$ws.Range("L58").Value = 8700
$ws.Range("M58").Value = -4264.0002
$ws.Range("N58").Value = -9000

$ws.Range("H69").Value = 5944

$ws.Range("H72").Value = 5944

$ws.Range("H116").Value = 44531.934
$ws.Range("I116").Value = 67553.336
$ws.Range("J116").Value = 9999.833000000001
$ws.Range("K116").Value = 67553.336
$ws.Range("L116").Value = 9999.833000000001
$ws.Range("M116").Value = -64111.336
$ws.Range("N116").Value = -16883.833

$ws.Range("H125").Value = 2208.8
$ws.Range("I125").Value = 1722.25
$ws.Range("J125").Value = 2533.1667
$ws.Range("K125").Value = 15500.25
$ws.Range("L125").Value = 22798.5003
$ws.Range("M125").Value = -13040.25
$ws.Range("N125").Value = -27718.5003

$ws.Range("H135").Value = 2150
$ws.Range("I135").Value = 1200
$ws.Range("J135").Value = 2466.6667
$ws.Range("K135").Value = 10800
$ws.Range("L135").Value = 22200.0003
$ws.Range("M135").Value = -8265
$ws.Range("N135").Value = -27270.0003

$ws.Range("H138").Value = 22226.574
$ws.Range("I138").Value = 43419.668
$ws.Range("J138").Value = 5272.1
$ws.Range("K138").Value = 130259.004
$ws.Range("L138").Value = 15816.3
$ws.Range("M138").Value = -125119.004
$ws.Range("N138").Value = -26096.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4761.243
$ws.Range("I32").Value = 4139.9697
$ws.Range("K32").Value = 4139.9697
$ws.Range("M32").Value = -3852.9697

$ws.Range("H74").Value = 7035.1562
$ws.Range("I74").Value = 1210.6818
$ws.Range("J74").Value = 19849
$ws.Range("K74").Value = 1210.6818
$ws.Range("L74").Value = 19849
$ws.Range("M74").Value = -336.6818000000001
$ws.Range("N74").Value = -21597

$ws.Range("H77").Value = 7035.1562
$ws.Range("I77").Value = 1210.6818
$ws.Range("J77").Value = 19849
$ws.Range("K77").Value = 6053.409000000001
$ws.Range("L77").Value = 99245
$ws.Range("M77").Value = -1685.409000000001
$ws.Range("N77").Value = -107981

$ws.Range("H122").Value = 2426.6365
$ws.Range("I122").Value = 2545.3333
$ws.Range("J122").Value = 2284.2
$ws.Range("K122").Value = 7635.999899999999
$ws.Range("L122").Value = 6852.599999999999
$ws.Range("M122").Value = -5185.999899999999
$ws.Range("N122").Value = -11752.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1688.175
$ws.Range("I107").Value = 1261.1
$ws.Range("K107").Value = 1261.1
$ws.Range("M107").Value = 658.9000000000001

$ws.Range("H134").Value = 1721.3864
$ws.Range("I134").Value = 1714.907
$ws.Range("K134").Value = 5144.721
$ws.Range("M134").Value = -2609.721

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 40833
$ws.Range("J50").Value = 51249.5
$ws.Range("L50").Value = 51249.5
$ws.Range("N50").Value = -52499.5

$ws.Range("H58").Value = 2656.2979
$ws.Range("I58").Value = 2246.889
$ws.Range("K58").Value = 2246.889
$ws.Range("M58").Value = -2043.889

$ws.Range("H60").Value = 8098.4287

$ws.Range("H99").Value = 3039.4211
$ws.Range("J99").Value = 3610.5
$ws.Range("L99").Value = 3610.5
$ws.Range("N99").Value = -6606.5

$ws.Range("H126").Value = 3039.4211
$ws.Range("J126").Value = 3610.5
$ws.Range("L126").Value = 10831.5
$ws.Range("N126").Value = -15771.5

$ws.Range("H132").Value = 106084.86
$ws.Range("I132").Value = 118708.375
$ws.Range("J132").Value = 15616.333
$ws.Range("K132").Value = 356125.125
$ws.Range("L132").Value = 46848.999
$ws.Range("M132").Value = -353595.125
$ws.Range("N132").Value = -51908.999

$ws.Range("H134").Value = 24490.191
$ws.Range("I134").Value = 17243.5
$ws.Range("K134").Value = 51730.5
$ws.Range("M134").Value = -49195.5

$ws.Range("H136").Value = 2656.2979
$ws.Range("I136").Value = 2246.889
$ws.Range("K136").Value = 6740.667
$ws.Range("M136").Value = -4190.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 528.6923
$ws.Range("J12").Value = 731.75
$ws.Range("L12").Value = 2195.25
$ws.Range("N12").Value = -2541.25

$ws.Range("H132").Value = 1135.7241

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H46").Value = 39542.25
$ws.Range("J46").Value = 39542.25
$ws.Range("L46").Value = 39542.25
$ws.Range("N46").Value = -39854.25

$ws.Range("H70").Value = 10620.8125
$ws.Range("I70").Value = 5540.4443
$ws.Range("J70").Value = 17152.715
$ws.Range("K70").Value = 5540.4443
$ws.Range("L70").Value = 17152.715
$ws.Range("M70").Value = -5270.4443
$ws.Range("N70").Value = -17692.715

$ws.Range("H73").Value = 10620.8125
$ws.Range("I73").Value = 5540.4443
$ws.Range("J73").Value = 17152.715
$ws.Range("K73").Value = 5540.4443
$ws.Range("L73").Value = 17152.715
$ws.Range("M73").Value = -4604.4443
$ws.Range("N73").Value = -19024.715

$ws.Range("H80").Value = 3403.2222
$ws.Range("I80").Value = 2949.6667
$ws.Range("J80").Value = 3630
$ws.Range("K80").Value = 2949.6667
$ws.Range("L80").Value = 3630
$ws.Range("M80").Value = -1951.6667
$ws.Range("N80").Value = -5626

$ws.Range("H83").Value = 3403.2222
$ws.Range("I83").Value = 2949.6667
$ws.Range("J83").Value = 3630
$ws.Range("K83").Value = 14748.3335
$ws.Range("L83").Value = 18150
$ws.Range("M83").Value = -9756.333500000001
$ws.Range("N83").Value = -28134

$ws.Range("H135").Value = 52777.445
$ws.Range("J135").Value = 52777.445
$ws.Range("L135").Value = 52777.445
$ws.Range("N135").Value = -62917.445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1199.8572
$ws.Range("I46").Value = 1183.1666
$ws.Range("K46").Value = 1183.1666
$ws.Range("M46").Value = -995.1666

$ws.Range("H136").Value = 42288.76
$ws.Range("I136").Value = 49477.094
$ws.Range("K136").Value = 148431.282
$ws.Range("M136").Value = -145881.282

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 49824.5
$ws.Range("J4").Value = 2239.75
$ws.Range("L4").Value = 2239.75
$ws.Range("N4").Value = -2465.75

$ws.Range("H81").Value = 9210.950000000001
$ws.Range("I81").Value = 25804.75
$ws.Range("J81").Value = 5062.5
$ws.Range("K81").Value = 51609.5
$ws.Range("L81").Value = 10125
$ws.Range("M81").Value = -50548.5
$ws.Range("N81").Value = -12247

$ws.Range("H84").Value = 9210.950000000001
$ws.Range("I84").Value = 25804.75
$ws.Range("J84").Value = 5062.5
$ws.Range("K84").Value = 258047.5
$ws.Range("L84").Value = 50625
$ws.Range("M84").Value = -252743.5
$ws.Range("N84").Value = -61233

$ws.Range("H132").Value = 1006.5625
$ws.Range("I132").Value = 950.2414
$ws.Range("J132").Value = 1551
$ws.Range("K132").Value = 2850.7242
$ws.Range("L132").Value = 4653
$ws.Range("M132").Value = -320.7242000000001
$ws.Range("N132").Value = -9713

$ws.Range("H136").Value = 2001.3636
$ws.Range("I136").Value = 1772.6177
$ws.Range("K136").Value = 5317.8531
$ws.Range("M136").Value = -2767.8531
